$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Resource" column (K) with a header and values for the 3 data rows.
$ws.Range("K1").Value = "Resource"
$ws.Range("K2").Value = "player_dahong"
$ws.Range("K3").Value = "player_huanying"
$ws.Range("K4").Value = "player_micai"

# Keep the selection consistent with the post-edit state shown in the diff (K4 selected).
$ws.Range("K4").Select()
